$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (unstyled data cell) used to restore format after forcing text entry
$plainStyle = $ws.Range("D6").Style

$ws.Range("D2").Value = '29.834.97'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '1.639.19'
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = "'215.32"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("E6").Value = '  -0.58%  '

$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").Value = "'28.86"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = '  -3.52%  '

$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("D11").Value = "'0.0899"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("D12").Value = '1.873.32'
$ws.Range("E12").Value = '  +0.77%  '

$ws.Range("D13").Value = '1.638.08'
$ws.Range("E13").Value = '  +0.66%  '

$ws.Range("D14").Value = "'0.593"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  +3.72%  '

$ws.Range("D15").Value = "'9.58"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  +8.06%  '

$ws.Range("D16").Value = "'3.90"
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = '  -0.30%  '

$ws.Range("D17").Value = '29.836.03'
$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").Value = "'64.29"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("D19").Value = "'237.71"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  -2.91%  '

$ws.Range("D20").Value = '0.0₃0703'
$ws.Range("E20").Value = '  -0.59%  '

$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").Value = "'9.92"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  +2.80%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").Value = "'2.18"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  +2.44%  '

$ws.Range("D25").Value = "'157.51"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  -0.72%  '

$ws.Range("D27").Value = "'0.110"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("D28").Value = "'6.64"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  +0.39%  '

$ws.Range("E29").Value = '  +0.29%  '

$ws.Range("D30").Value = "'0.0495"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  +1.07%  '

$ws.Range("E31").Value = '  -1.07%  '

$ws.Range("E32").Value = '  +0.68%  '

$ws.Range("E33").Value = '  -1.42%  '

$ws.Range("D34").Value = '1.420.21'
$ws.Range("E34").Value = '  -0.38%  '

$ws.Range("E35").Value = '  +2.09%  '

$ws.Range("E36").Value = '  -1.73%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.0173"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = "'2.66"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  -7.15%  '

$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").Value = "'76.63"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  +10.27%  '

$ws.Range("D41").Value = "'0.567"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  +1.36%  '

$ws.Range("D42").Value = "'0.0503"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").Value = "'0.833"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = '  -0.13%  '

$ws.Range("E44").Value = '  -2.69%  '

$ws.Range("E45").Value = '  +0.39%  '

$ws.Range("E46").Value = '  -2.22%  '

$ws.Range("D47").Value = "'50.22"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  -8.46%  '

$ws.Range("D48").Value = '1.781.67'
$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("E49").Value = '  -1.61%  '

$ws.Range("D50").Value = "'93.58"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  +5.40%  '

$ws.Range("D51").Value = '0.0₆0111'
$ws.Range("E51").Value = '  +2.69%  '
